$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.066.27'
$ws.Range('E2').Value = '  +0.20%  '

$ws.Range('D3').Value = '2.959.35'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '379.93'
$ws.Range('E5').Value = '  +0.93%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.51'
$ws.Range('E6').Value = '  +0.21%  '

$ws.Range('E7').Value = '  +1.83%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.587'
$ws.Range('E9').Value = '  +0.68%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.46'
$ws.Range('E10').Value = '  -0.24%  '

$ws.Range('E11').Value = '  -0.49%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0852'
$ws.Range('E12').Value = '  +2.05%  '

$ws.Range('D13').Value = '3.428.25'
$ws.Range('E13').Value = '  +1.06%  '

$ws.Range('B14').Value = 'Uniswap'
$ws.Range('C14').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '12.42'
$ws.Range('E14').Value = '  +74.42%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '18.41'
$ws.Range('E15').Value = '  +2.67%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.74'
$ws.Range('E16').Value = '  +5.57%  '

$ws.Range('D17').Value = '2.952.84'
$ws.Range('E17').Value = '  +0.94%  '

$ws.Range('E18').Value = '  +4.04%  '

$ws.Range('D19').Value = '51.115.47'
$ws.Range('E19').Value = '  +0.48%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.06'
$ws.Range('E20').Value = '  -3.15%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.41'
$ws.Range('E21').Value = '  -0.61%  '

$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('E23').Value = '  +17.57%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.60'
$ws.Range('E24').Value = '  +2.92%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.69'
$ws.Range('E25').Value = '  +2.39%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.98'
$ws.Range('E26').Value = '  -2.85%  '

$ws.Range('E27').Value = '  -0.06%  '

$ws.Range('E28').Value = '  -0.80%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '25.87'
$ws.Range('E29').Value = '  +1.01%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -9.46%  '

$ws.Range('E31').Value = '  -4.36%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.54'
$ws.Range('E32').Value = '  +7.52%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.29'
$ws.Range('E33').Value = '  +1.06%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '50.94'
$ws.Range('E34').Value = '  +0.68%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.05'
$ws.Range('E35').Value = '  +1.57%  '

$ws.Range('E36').Value = '  -4.65%  '

$ws.Range('E37').Value = '  +0.00%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.24'
$ws.Range('E38').Value = '  +8.99%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '16.78'
$ws.Range('E39').Value = '  +2.37%  '

$ws.Range('E40').Value = '  +2.07%  '

$ws.Range('E41').Value = '  +2.79%  '

$ws.Range('E42').Value = '  -2.35%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '121.67'
$ws.Range('E43').Value = '  +0.37%  '

$ws.Range('E44').Value = '  +11.80%  '

$ws.Range('E45').Value = '  +2.64%  '

$ws.Range('E46').Value = '  -1.03%  '

$ws.Range('D47').Value = '2.055.54'
$ws.Range('E47').Value = '  +2.62%  '

$ws.Range('E48').Value = '  -1.00%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.257'
$ws.Range('E49').Value = '  -5.14%  '

$ws.Range('E50').Value = '  -5.76%  '

$ws.Range('E51').Value = '  +6.87%  '

